# "new links for examples"
# Update the "Ejemplo de ..." (Example) links in the "table" sheet so that
# each recommendation row now points to a numbered (1.1 / 1.2) example link
# instead of a single shared example per category.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table")

$ws.Range("E2").Value  = "1. 1 Ejemplo de base de datos"
$ws.Range("E3").Value  = "1. 1 Ejemplo de base de datos"
$ws.Range("E4").Value  = "1.2 Ejemplo de base de datos"
$ws.Range("E5").Value  = "1.2 Ejemplo de base de datos"

$ws.Range("E6").Value  = "1.1 Ejemplo de cuestionario"
$ws.Range("E7").Value  = "1.2 Ejemplo de cuestionario"

$ws.Range("E8").Value  = "1.1 Ejemplo de libro de códigos"
$ws.Range("E9").Value  = "1.2 Ejemplo de libro de códigos"

$ws.Range("E10").Value = "1.1 Ejemplo de ficha técnica"
$ws.Range("E11").Value = "1.1 Ejemplo de ficha técnica"
$ws.Range("E12").Value = "1.2 Ejemplo de ficha técnica"
$ws.Range("E13").Value = "1.2 Ejemplo de ficha técnica"

# Selection / scroll position changed in the saved file (cosmetic, but
# reproduce it so the sheetView matches the new session state).
$ws.Activate()
$ws.Range("L13").Select()

# The "pre" sheet's merged-cell list was re-ordered on save (A/B columns
# moved ahead of D/E). Re-apply the merges in that order to reproduce it.
$ws2 = $wb.Worksheets.Item("pre")
$mergeRanges = @(
    "A2:A5","A6:A7","A8:A9","A10:A13",
    "B2:B3","B4:B5","B6:B7","B8:B9","B10:B13",
    "D2:D5","D6:D13",
    "E2:E5","E6:E7","E8:E9","E10:E13"
)
foreach ($r in $mergeRanges) {
    $ws2.Range($r).UnMerge()
}
foreach ($r in $mergeRanges) {
    $ws2.Range($r).Merge()
}
